$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Liste des tâches")

# Mark the already-completed tasks ("Fait") for the two user stories that
# have been finished: "Enregistrement et connexion utilisateurs" (rows 2-7)
# and "Connexion et configuration du serveur Core" (rows 47-49).
$doneRows = @(2, 3, 4, 6, 7, 47, 48, 49)
foreach ($r in $doneRows) {
    $cell = $ws.Range("D$r")
    $cell.Value = "Fait"
    $cell.Font.Color = 255
}

# This task (ConnectUserBean / authentication) is not done yet, but still
# picks up the same (empty) styled cell.
$ws.Range("D5").Font.Color = 255

# Reflect the page setup used when this sheet was last printed/reviewed.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Make "Liste des tâches" the active sheet/tab, scrolled back to the top,
# with the last worked-on cell selected.
$ws.Select()
$ws.Range("G48").Select()
